$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prueba")

# "harina" -> "harinita"
$ws.Range("A2").Value = "harinita"

# Quantity for row 2: 2.0 -> 1.0
$ws.Range("B2").Value = 1.0

# New row 11: "test" / 12.0
$ws.Range("A11").Value = "test"
$ws.Range("B11").Value = 12.0
